$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to be treated as text so numeric-looking
# strings (e.g. "1.00", "236.90") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.975.38"
$ws.Range("E2").Value = "  -1.37%  "

# Row 3
$ws.Range("D3").Value = "2.224.47"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
$ws.Range("D5").Value = "299.27"
$ws.Range("E5").Value = "  -3.19%  "

# Row 6
$ws.Range("D6").Value = "89.81"
$ws.Range("E6").Value = "  -6.16%  "

# Row 7
$ws.Range("D7").Value = "0.557"
$ws.Range("E7").Value = "  -2.99%  "

# Row 8
$ws.Range("E8").Value = "  -0.55%  "

# Row 10
$ws.Range("D10").Value = "32.63"
$ws.Range("E10").Value = "  -7.44%  "

# Row 11
$ws.Range("D11").Value = "0.0777"
$ws.Range("E11").Value = "  -4.14%  "

# Row 12
$ws.Range("D12").Value = "6.95"
$ws.Range("E12").Value = "  -4.88%  "

# Row 13
$ws.Range("E13").Value = "  -0.97%  "

# Row 14
$ws.Range("D14").Value = "2.564.01"
$ws.Range("E14").Value = "  -1.14%  "

# Row 15
$ws.Range("D15").Value = "2.210.13"
$ws.Range("E15").Value = "  -1.96%  "

# Row 16
$ws.Range("D16").Value = "13.55"
$ws.Range("E16").Value = "  -1.00%  "

# Row 17
$ws.Range("D17").Value = "0.776"
$ws.Range("E17").Value = "  -7.96%  "

# Row 18
$ws.Range("D18").Value = "43.901.50"
$ws.Range("E18").Value = "  -0.83%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0902"
$ws.Range("E19").Value = "  -6.72%  "

# Row 20
$ws.Range("D20").Value = "5.91"

# Row 21
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  -9.59%  "

# Row 22
$ws.Range("D22").Value = "64.67"
$ws.Range("E22").Value = "  -1.95%  "

# Row 23
$ws.Range("D23").Value = "236.90"
$ws.Range("E23").Value = "  -1.64%  "

# Row 24
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -6.53%  "

# Row 25
$ws.Range("E25").Value = "  +0.51%  "

# Row 26
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  -7.56%  "

# Row 27
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").Value = "38.22"
$ws.Range("E27").Value = "  +1.18%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.17"
$ws.Range("E28").Value = "  -3.04%  "

# Row 29
$ws.Range("E29").Value = "  -6.04%  "

# Row 30
$ws.Range("D30").Value = "19.27"
$ws.Range("E30").Value = "  -4.50%  "

# Row 31
$ws.Range("D31").Value = "149.83"
$ws.Range("E31").Value = "  -1.91%  "

# Row 32
$ws.Range("E32").Value = "  -10.51%  "

# Row 33
$ws.Range("E33").Value = "  -6.79%  "

# Row 34
$ws.Range("E34").Value = "  -4.81%  "

# Row 35
$ws.Range("E35").Value = "  -4.03%  "

# Row 36
$ws.Range("D36").Value = "2.80"
$ws.Range("E36").Value = "  -11.92%  "

# Row 37
$ws.Range("E37").Value = "  -7.85%  "

# Row 38
$ws.Range("D38").Value = "1.68"
$ws.Range("E38").Value = "  -7.26%  "

# Row 39
$ws.Range("E39").Value = "  -1.26%  "

# Row 40
$ws.Range("D40").Value = "3.17"
$ws.Range("E40").Value = "  -8.37%  "

# Row 41
$ws.Range("D41").Value = "3.53"
$ws.Range("E41").Value = "  -7.72%  "

# Row 42
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.88%  "

# Row 43
$ws.Range("D43").Value = "13.01"
$ws.Range("E43").Value = "  -10.94%  "

# Row 44
$ws.Range("D44").Value = "1.806.25"
$ws.Range("E44").Value = "  +3.13%  "

# Row 45
$ws.Range("D45").Value = "1.79"
$ws.Range("E45").Value = "  +12.50%  "

# Row 46
$ws.Range("D46").Value = "0.178"

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "14.18"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "73.31"
$ws.Range("E48").Value = "  -10.10%  "

# Row 49
$ws.Range("D49").Value = "93.79"
$ws.Range("E49").Value = "  -6.31%  "

# Row 50
$ws.Range("D50").Value = "66.56"
$ws.Range("E50").Value = "  -6.75%  "

# Row 51
$ws.Range("D51").Value = "2.445.43"
$ws.Range("E51").Value = "  -1.16%  "
